$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Text
    )
    $range = $ws.Range($Cell)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.Style = "Normal"
}

Set-TextValue "D2" "26.894.41"
Set-TextValue "E2" "  +1.85%  "
Set-TextValue "D3" "1.726.73"
Set-TextValue "E3" "  +0.27%  "
Set-TextValue "D4" "0.9973"
Set-TextValue "E4" "  -0.27%  "
Set-TextValue "D5" "242.22"
Set-TextValue "E5" "  -0.16%  "
Set-TextValue "D6" "0.9977"
Set-TextValue "E6" "  -0.23%  "
Set-TextValue "D7" "0.4895"
Set-TextValue "E7" "  -0.40%  "
Set-TextValue "D8" "0.2592"
Set-TextValue "E8" "  -0.77%  "
Set-TextValue "D9" "0.06212"
Set-TextValue "E9" "  +0.33%  "
Set-TextValue "D10" "1.730.95"
Set-TextValue "E10" "  +0.54%  "
Set-TextValue "D11" "15.98"
Set-TextValue "E11" "  +3.27%  "
Set-TextValue "D12" "0.06905"
Set-TextValue "D13" "0.6084"
Set-TextValue "E13" "  +1.67%  "
Set-TextValue "D14" "4.487"
Set-TextValue "E14" "  -1.62%  "
Set-TextValue "D15" "77.26"
Set-TextValue "E15" "  +0.10%  "
Set-TextValue "D16" "0.9982"
Set-TextValue "E16" "  -0.22%  "
Set-TextValue "D17" "26.647.44"
Set-TextValue "E17" "  +0.92%  "
Set-TextValue "E18" "  -0.29%  "
Set-TextValue "E19" "  +0.50%  "
Set-TextValue "E20" "  +0.86%  "
Set-TextValue "D21" "1.953.76"
Set-TextValue "E21" "  +0.52%  "
Set-TextValue "D22" "4.425"
Set-TextValue "E22" "  -1.17%  "
Set-TextValue "D23" "8.567"
Set-TextValue "E23" "  -0.10%  "
Set-TextValue "D24" "5.099"
Set-TextValue "E24" "  -1.02%  "
Set-TextValue "D25" "138.48"
Set-TextValue "E25" "  +0.89%  "
Set-TextValue "D26" "15.33"
Set-TextValue "E26" "  +0.76%  "
Set-TextValue "D27" "1.782"
Set-TextValue "E27" "  +4.73%  "
$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D28" "106.38"
Set-TextValue "E28" "  -0.56%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D29" "1.380"
Set-TextValue "E29" "  -1.22%  "
Set-TextValue "D30" "3.946"
Set-TextValue "E30" "  +0.22%  "
Set-TextValue "D31" "0.08000"
Set-TextValue "E31" "  +0.62%  "
Set-TextValue "E32" "  +0.59%  "
Set-TextValue "D33" "0.04524"
Set-TextValue "E33" "  -0.43%  "
Set-TextValue "D34" "2.596"
Set-TextValue "E34" "  -0.22%  "
Set-TextValue "E35" "  +1.60%  "
Set-TextValue "D36" "0.6247"
Set-TextValue "E36" "  +0.21%  "
Set-TextValue "D37" "0.9383"
Set-TextValue "E37" "  +1.61%  "
Set-TextValue "D38" "2.054"
Set-TextValue "E38" "  +5.62%  "
Set-TextValue "D39" "2.451"
Set-TextValue "E39" "  +2.48%  "
Set-TextValue "D40" "0.9977"
Set-TextValue "E40" "  -0.20%  "
Set-TextValue "D41" "0.01501"
Set-TextValue "E41" "  +1.38%  "
Set-TextValue "D42" "5.665"
Set-TextValue "E42" "  +6.04%  "
Set-TextValue "D43" "99.54"
Set-TextValue "E43" "  -0.42%  "
Set-TextValue "D44" "0.3851"
Set-TextValue "E44" "  +0.52%  "
Set-TextValue "D45" "6.873"
Set-TextValue "E45" "  +2.52%  "
Set-TextValue "D46" "0.1161"
Set-TextValue "E46" "  -0.11%  "
Set-TextValue "D47" "0.05399"
Set-TextValue "E47" "  +0.67%  "
Set-TextValue "D48" "7.913"
Set-TextValue "E48" "  +3.14%  "
Set-TextValue "D49" "30.14"
Set-TextValue "E49" "  +0.24%  "
Set-TextValue "D50" "51.65"
Set-TextValue "E50" "  +1.64%  "
Set-TextValue "D51" "1.234"
Set-TextValue "E51" "  +0.10%  "
